$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 88, shifting all
# subsequent rows (old 88..204) down by one (new 89..205).
$ws.Rows("88:88").Insert()

# Populate the newly inserted row 88 with the new weekly record.
$ws.Cells.Item(88, 1).Value = 5
$ws.Cells.Item(88, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(88, 3).Value = "Maule"
$ws.Cells.Item(88, 4).Value = 44930
$ws.Cells.Item(88, 5).Value = 7
$ws.Cells.Item(88, 6).Value = 100112031
$ws.Cells.Item(88, 7).Value = "Poroto verde"
$ws.Cells.Item(88, 8).Value = "Sin especificar"
$ws.Cells.Item(88, 9).Value = "Primera"
$ws.Cells.Item(88, 10).Value = 200
$ws.Cells.Item(88, 11).Value = 30000
$ws.Cells.Item(88, 12).Value = 30000
$ws.Cells.Item(88, 13).Value = 30000
$ws.Cells.Item(88, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(88, 15).Value = "Región del Maule"
$ws.Cells.Item(88, 16).Value = 1200
$ws.Cells.Item(88, 17).Value = 25
$ws.Cells.Item(88, 18).Value = "Hortaliza"
